# Append " (Changed main)" after the existing sentence in the first
# paragraph, as three additional runs:
#   <w:r><w:t xml:space="preserve"> (</w:t></w:r>
#   <w:r><w:t>Changed main</w:t></w:r>
#   <w:r><w:t>)</w:t></w:r>
#
# A plain Range.InsertAfter() would just grow the existing run's text
# (adjacent runs that share formatting get coalesced into one run), so
# track changes is toggled on for the duration of the inserts - each
# InsertAfter then lands as its own <w:r> inside a <w:ins> - and every
# resulting insertion revision is accepted individually (not via
# Revisions.AcceptAll, which forces a full-document relayout). That
# leaves three separate sibling runs behind, matching the diff, without
# perturbing any other part of the document.

$d = $word.ActiveDocument

$d.TrackRevisions = $true

$para = $d.Paragraphs(1).Range
$insertionPoint = $para.End - 1

$r1 = $d.Range($insertionPoint, $insertionPoint)
$r1.InsertAfter(" (")

$insertionPoint = $insertionPoint + 2
$r2 = $d.Range($insertionPoint, $insertionPoint)
$r2.InsertAfter("Changed main")

$insertionPoint = $insertionPoint + 12
$r3 = $d.Range($insertionPoint, $insertionPoint)
$r3.InsertAfter(")")

$d.TrackRevisions = $false

for ($i = $d.Revisions.Count; $i -ge 1; $i--) {
    $d.Revisions($i).Accept()
}
